$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 621.9167
$ws.Range("I6").Value = 273.66666
$ws.Range("J6").Value = 1666.6666
$ws.Range("K6").Value = 820.9999799999999
$ws.Range("L6").Value = 4999.9998
$ws.Range("M6").Value = -708.9999799999999
$ws.Range("N6").Value = -5223.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3788
$ws.Range("I74").Value = 3550
$ws.Range("J74").Value = 4502
$ws.Range("K74").Value = 3550
$ws.Range("L74").Value = 4502
$ws.Range("M74").Value = -2614
$ws.Range("N74").Value = -6374

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 18521848
$ws.Range("I76").Value = 111111110
$ws.Range("K76").Value = 111111110
$ws.Range("M76").Value = -111110795

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3788
$ws.Range("I77").Value = 3550
$ws.Range("J77").Value = 4502
$ws.Range("K77").Value = 17750
$ws.Range("L77").Value = 22510
$ws.Range("M77").Value = -13070
$ws.Range("N77").Value = -31870

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 18521848
$ws.Range("I79").Value = 111111110
$ws.Range("K79").Value = 111111110
$ws.Range("M79").Value = -111110018

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2158210.5
$ws.Range("J116").Value = 2502.5
$ws.Range("L116").Value = 2502.5
$ws.Range("N116").Value = -9386.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 9781633
$ws.Range("I138").Value = 2648598.5
$ws.Range("J138").Value = 13160439
$ws.Range("K138").Value = 7945795.5
$ws.Range("L138").Value = 39481317
$ws.Range("M138").Value = -7940655.5
$ws.Range("N138").Value = -39491597

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 38126.5
$ws.Range("I6").Value = 38126.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 38126.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -37953.5
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6285.3667
$ws.Range("I63").Value = 8545.421
$ws.Range("J63").Value = 2381.6365
$ws.Range("K63").Value = 8545.421
$ws.Range("L63").Value = 2381.6365
$ws.Range("M63").Value = -7859.421
$ws.Range("N63").Value = -3753.6365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6285.3667
$ws.Range("I66").Value = 8545.421
$ws.Range("J66").Value = 2381.6365
$ws.Range("K66").Value = 42727.105
$ws.Range("L66").Value = 11908.1825
$ws.Range("M66").Value = -39295.105
$ws.Range("N66").Value = -18772.1825

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 261.66666
$ws.Range("I22").Value = 213.07692
$ws.Range("K22").Value = 213.07692
$ws.Range("M22").Value = -40.07692

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 280926.97
$ws.Range("I105").Value = 3023.913
$ws.Range("J105").Value = 772601.6
$ws.Range("K105").Value = 3023.913
$ws.Range("L105").Value = 772601.6
$ws.Range("M105").Value = -1276.913
$ws.Range("N105").Value = -776095.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11065.214
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 11065.214
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 11065.214
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -11655.214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11065.214
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 11065.214
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 11065.214
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -11469.214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 13412.238
$ws.Range("I62").Value = 24025.5
$ws.Range("J62").Value = 3763.818
$ws.Range("K62").Value = 24025.5
$ws.Range("L62").Value = 3763.818
$ws.Range("M62").Value = -23401.5
$ws.Range("N62").Value = -5011.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 13412.238
$ws.Range("I65").Value = 24025.5
$ws.Range("J65").Value = 3763.818
$ws.Range("K65").Value = 120127.5
$ws.Range("L65").Value = 18819.09
$ws.Range("M65").Value = -117007.5
$ws.Range("N65").Value = -25059.09

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3742.0645
$ws.Range("I132").Value = 3635.95
$ws.Range("J132").Value = 3935
$ws.Range("K132").Value = 10907.85
$ws.Range("L132").Value = 11805
$ws.Range("M132").Value = -8377.849999999999
$ws.Range("N132").Value = -16865

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5377661.5
$ws.Range("J131").Value = 5849351.5
$ws.Range("L131").Value = 17548054.5
$ws.Range("N131").Value = -17558134.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 50000
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6742.4
$ws.Range("I70").Value = 7295.273
$ws.Range("J70").Value = 6066.6665
$ws.Range("K70").Value = 7295.273
$ws.Range("L70").Value = 6066.6665
$ws.Range("M70").Value = -7025.273
$ws.Range("N70").Value = -6606.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6742.4
$ws.Range("I73").Value = 7295.273
$ws.Range("J73").Value = 6066.6665
$ws.Range("K73").Value = 7295.273
$ws.Range("L73").Value = 6066.6665
$ws.Range("M73").Value = -6359.273
$ws.Range("N73").Value = -7938.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 31252874
$ws.Range("I80").Value = 2991.5386
$ws.Range("J80").Value = 166669040
$ws.Range("K80").Value = 2991.5386
$ws.Range("L80").Value = 166669040
$ws.Range("M80").Value = -1993.5386
$ws.Range("N80").Value = -166671036

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 31252874
$ws.Range("I83").Value = 2991.5386
$ws.Range("J83").Value = 166669040
$ws.Range("K83").Value = 14957.693
$ws.Range("L83").Value = 833345200
$ws.Range("M83").Value = -9965.693
$ws.Range("N83").Value = -833355184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 13590.818
$ws.Range("J123").Value = 14314.9
$ws.Range("L123").Value = 14314.9
$ws.Range("N123").Value = -19214.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1702.975
$ws.Range("I132").Value = 1673.4375
$ws.Range("J132").Value = 1821.125
$ws.Range("K132").Value = 5020.3125
$ws.Range("L132").Value = 5463.375
$ws.Range("M132").Value = -2490.3125
$ws.Range("N132").Value = -10523.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 5333.3335
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12227.777
$ws.Range("I22").Value = 1233.3334
$ws.Range("J22").Value = 17725
$ws.Range("K22").Value = 1233.3334
$ws.Range("L22").Value = 17725
$ws.Range("M22").Value = -938.3334
$ws.Range("N22").Value = -18315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 12227.777
$ws.Range("I27").Value = 1233.3334
$ws.Range("J27").Value = 17725
$ws.Range("K27").Value = 1233.3334
$ws.Range("L27").Value = 17725
$ws.Range("M27").Value = -1126.3334
$ws.Range("N27").Value = -17939

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2838.75
$ws.Range("I40").Value = 1482.5
$ws.Range("K40").Value = 1482.5
$ws.Range("M40").Value = -1346.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 500.5
$ws.Range("J55").Value = 470.75
$ws.Range("L55").Value = 470.75
$ws.Range("N55").Value = -816.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2399.077
$ws.Range("I132").Value = 1987.9642
$ws.Range("K132").Value = 5963.892599999999
$ws.Range("M132").Value = -3433.892599999999
